# Unlock Conversation Effect and Unlock Note Effect Update
# Insert two new option rows (308, 309) before the existing terminal rows
# (998 "Skip"/999 "Generic end"), pushing those down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 35 (shifts old rows 35-36 down to 37-38)
$ws.Rows.Item(35).Insert()
$ws.Rows.Item(35).Insert()

# New row 35: Option 308 - Unlock new note
$ws.Cells.Item(35, 1).Value = 308
$ws.Cells.Item(35, 2).Value = "解锁新笔记"
$ws.Cells.Item(35, 3).Value = "解锁新笔记"
$ws.Cells.Item(35, 4).Value = "[[105:101:3]]"

# New row 36: Option 309 - Unlock new conversation
$ws.Cells.Item(36, 1).Value = 309
$ws.Cells.Item(36, 2).Value = "解锁新对话"
$ws.Cells.Item(36, 3).Value = "解锁新对话"
$ws.Cells.Item(36, 4).Value = "[[106:1001:1001]]"

# Match the author's final selection state after editing
$ws.Range("E35").Select()
